# Auto-generated: apply cryptos list update (Sun Sep  8 19:19:57 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, [string]$val, [bool]$forceText = $false)
    if ($forceText) {
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $val
        $ws.Range($addr).Style = "Normal"
    } else {
        $ws.Range($addr).Value = $val
    }
}

# Row 2
Set-CellText $ws 'D2' '54.313.80' $false
Set-CellText $ws 'E2' '  +0.54%  ' $false

# Row 3
Set-CellText $ws 'D3' '2.265.42' $false
Set-CellText $ws 'E3' '  -0.89%  ' $false

# Row 4
Set-CellText $ws 'E4' '  +0.10%  ' $false

# Row 5
Set-CellText $ws 'D5' '496.28' $true
Set-CellText $ws 'E5' '  +0.14%  ' $false

# Row 6
Set-CellText $ws 'D6' '128.92' $true
Set-CellText $ws 'E6' '  +0.86%  ' $false

# Row 7
Set-CellText $ws 'E7' '  +0.08%  ' $false

# Row 8
Set-CellText $ws 'E8' '  -0.70%  ' $false

# Row 9
Set-CellText $ws 'E9' '  +0.64%  ' $false

# Row 10
Set-CellText $ws 'E10' '  +1.00%  ' $false

# Row 11
Set-CellText $ws 'E11' '  +3.50%  ' $false

# Row 12
Set-CellText $ws 'D12' '4.79' $true
Set-CellText $ws 'E12' '  +3.18%  ' $false

# Row 13
Set-CellText $ws 'D13' '23.02' $true
Set-CellText $ws 'E13' '  +5.56%  ' $false

# Row 14
Set-CellText $ws 'D14' '2.669.48' $false
Set-CellText $ws 'E14' '  +0.42%  ' $false

# Row 15
Set-CellText $ws 'D15' '54.295.74' $false
Set-CellText $ws 'E15' '  +0.55%  ' $false

# Row 16
Set-CellText $ws 'D16' '0.0000129' $true
Set-CellText $ws 'E16' '  +0.37%  ' $false

# Row 17
Set-CellText $ws 'D17' '2.267.49' $false
Set-CellText $ws 'E17' '  -1.17%  ' $false

# Row 18
Set-CellText $ws 'E18' '  +2.20%  ' $false

# Row 19
Set-CellText $ws 'D19' '4.13' $true
Set-CellText $ws 'E19' '  +1.15%  ' $false

# Row 20
Set-CellText $ws 'D20' '303.02' $true
Set-CellText $ws 'E20' '  +0.67%  ' $false

# Row 21
Set-CellText $ws 'D21' '6.33' $true
Set-CellText $ws 'E21' '  -1.63%  ' $false

# Row 22
Set-CellText $ws 'E22' '  +0.38%  ' $false

# Row 23
Set-CellText $ws 'D23' '60.88' $true
Set-CellText $ws 'E23' '  -2.58%  ' $false

# Row 24
Set-CellText $ws 'E24' '  -1.30%  ' $false

# Row 25
Set-CellText $ws 'E25' '  +1.00%  ' $false

# Row 26
Set-CellText $ws 'D26' '7.32' $true
Set-CellText $ws 'E26' '  +3.69%  ' $false

# Row 27
Set-CellText $ws 'D27' '171.52' $true
Set-CellText $ws 'E27' '  +1.46%  ' $false

# Row 28
Set-CellText $ws 'D28' '5.98' $true
Set-CellText $ws 'E28' '  +2.15%  ' $false

# Row 29
Set-CellText $ws 'E29' '  +0.01%  ' $false

# Row 30
Set-CellText $ws 'D30' '0.0₃0690' $false
Set-CellText $ws 'E30' '  +0.38%  ' $false

# Row 31
Set-CellText $ws 'E31' '  +1.39%  ' $false

# Row 32
Set-CellText $ws 'E32' '  -0.02%  ' $false

# Row 33
Set-CellText $ws 'D33' '17.82' $true
Set-CellText $ws 'E33' '  +0.81%  ' $false

# Row 34
Set-CellText $ws 'D34' '0.997' $true
Set-CellText $ws 'E34' '  +0.05%  ' $false

# Row 35
Set-CellText $ws 'D35' '0.936' $true
Set-CellText $ws 'E35' '  +8.01%  ' $false

# Row 36
Set-CellText $ws 'E36' '  +0.09%  ' $false

# Row 37
Set-CellText $ws 'E37' '  +0.14%  ' $false

# Row 38
Set-CellText $ws 'D38' '0.375' $true
Set-CellText $ws 'E38' '  +0.36%  ' $false

# Row 39
Set-CellText $ws 'E39' '  -0.24%  ' $false

# Row 40
Set-CellText $ws 'E40' '  +0.58%  ' $false

# Row 41
Set-CellText $ws 'B41' 'RenderToken' $false
Set-CellText $ws 'C41' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' $false
Set-CellText $ws 'D41' '4.81' $true
Set-CellText $ws 'E41' '  -2.50%  ' $false

# Row 42
Set-CellText $ws 'B42' 'Aave' $false
Set-CellText $ws 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' $false
Set-CellText $ws 'D42' '125.18' $true
Set-CellText $ws 'E42' '  -1.69%  ' $false

# Row 44
Set-CellText $ws 'E44' '  +0.77%  ' $false

# Row 45
Set-CellText $ws 'E45' '  +0.67%  ' $false

# Row 46
Set-CellText $ws 'D46' '241.96' $true
Set-CellText $ws 'E46' '  +1.44%  ' $false

# Row 47
Set-CellText $ws 'E47' '  +0.41%  ' $false

# Row 48
Set-CellText $ws 'E48' '  +1.35%  ' $false

# Row 49
Set-CellText $ws 'D49' '10.80' $true
Set-CellText $ws 'E49' '  +0.83%  ' $false

# Row 50
Set-CellText $ws 'D50' '16.10' $true
Set-CellText $ws 'E50' '  -0.91%  ' $false

# Row 51
Set-CellText $ws 'D51' '0.935' $true
Set-CellText $ws 'E51' '  -0.44%  ' $false

